# Append: 2025-10-29 06:36 JST
# Update the acquisition-timestamp column (A2:A10) on the active sheet
# from the previous run's value to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-29 06:36:01"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
